# Fruta / hortaliza, semanal
# Insert this week's new price observation as a new row 159 (most recent date
# first), shifting all the existing data rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row before row 159; everything currently at/after row 159
# (through the old last row 211) shifts down to 160..212.
$ws.Rows.Item(159).Insert()

# Populate the newly inserted row with the latest weekly record.
$ws.Cells.Item(159, 1).Value  = 4
$ws.Cells.Item(159, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(159, 3).Value  = "Los Lagos"
$ws.Cells.Item(159, 4).Value  = 44809
$ws.Cells.Item(159, 5).Value  = 10
$ws.Cells.Item(159, 6).Value  = 100112009
$ws.Cells.Item(159, 7).Value  = "Acelga"
$ws.Cells.Item(159, 8).Value  = "Sin especificar"
$ws.Cells.Item(159, 9).Value  = "Primera"
$ws.Cells.Item(159, 10).Value = 50
$ws.Cells.Item(159, 11).Value = 3000
$ws.Cells.Item(159, 12).Value = 3000
$ws.Cells.Item(159, 13).Value = 3000
$ws.Cells.Item(159, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(159, 15).Value = "Región del Maule"
$ws.Cells.Item(159, 16).Value = 500
$ws.Cells.Item(159, 17).Value = 6
$ws.Cells.Item(159, 18).Value = "Hortaliza"
